# Quick Start Tutorial - How to Compile: small update (#859)
#
# Removes the "After compilation... / Now you can add the agent...
# / Load the agent using the GUI... / (image)" block of paragraphs,
# and clears the "Figure 13: Adding an agent using the GUI." caption
# text while leaving its (now empty) paragraph in place.

$d = $word.ActiveDocument

# Locate the paragraph that starts the block to remove, and the
# paragraph right after the one containing the "...main class is the
# class that extends the..." text (i.e. the image paragraph), which
# marks the end of the block to remove.
$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "After compilation, you can move the entire package around*") {
        $startPara = $p
    }
    if ($t -like "*main class is the class that extends the*") {
        $endPara = $p.Next()
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}

# Re-locate the "Figure 13" caption paragraph now that earlier
# paragraphs/ranges have shifted, and blank out its text while
# keeping the (now empty) paragraph mark / formatting in place.
$figPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Figure 13: Adding an agent using the GUI.*") {
        $figPara = $p
    }
}

if ($figPara -ne $null) {
    $figRange = $figPara.Range
    $captionTextRange = $d.Range($figRange.Start, $figRange.End - 1)
    $captionTextRange.Delete()
}
